$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 2506534.94
$ws.Range("C9").Value = 357966.99
$ws.Range("D9").Value = 2864501.93
$ws.Range("E9").Value = 12.49665731591949
$ws.Range("F9").Value = 87.50334268408052
$ws.Range("G9").Value = -65.40427532219925
$ws.Range("H9").Value = -54.73544361548927
$ws.Range("I9").Value = -56.41511453061725
$ws.Range("J9").Value = 24802
$ws.Range("K9").Value = 1034
$ws.Range("L9").Value = 25836
